$d = $word.ActiveDocument

# Merge the two bold runs ("DOCX, DOC, PDF, HTML, XPS, R" and "TF and TXT",
# which were separated by a _GoBack bookmark) into a single run containing
# the combined text "DOCX, DOC, PDF, HTML, XPS, RTF and TXT". Replacing the
# whole found range also removes the bookmark that used to sit in the middle.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("DOCX, DOC, PDF, HTML, XPS, R" + "TF and TXT", $true, $false, $false, $false, $false, `
              $true, 1, $false, "DOCX, DOC, PDF, HTML, XPS, RTF and TXT", 2)
